# Update "想去人数" (number of people interested) figures for several events
# across the "展览", "演出" and "全部类型" sheets, per latest scrape run.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 816
$wsExhibition.Range("F5").Value = 982
$wsExhibition.Range("F6").Value = 2337

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 6

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 816
$wsAll.Range("F7").Value = 982
$wsAll.Range("F8").Value = 2337
$wsAll.Range("F9").Value = 6
